$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20 will become the new last row of the "Periodo Mora" detail table, so
# give it the same closing-border look the previous last row (28) had,
# before that row disappears.
$ws.Range("B28:J28").Copy()
$ws.Range("B20:J20").PasteSpecial(-4122)

# Drop the 8 oldest "Periodo Mora" rows (2106..2011); the signature block
# below shifts up automatically to close the gap.
$ws.Rows("21:28").Delete()

# Rewrite the remaining 5 detail rows with the new period data (ascending
# this time: 1909, 1910, 1911, 1912, 2001) and their "Valor Mora" amounts.
$ws.Range("E16").Value = "1909"
$ws.Range("F16").Value = 12146

$ws.Range("E17").Value = "1910"
$ws.Range("F17").Value = 33125

$ws.Range("E18").Value = "1911"
$ws.Range("F18").Value = 33125

$ws.Range("E19").Value = "1912"
$ws.Range("F19").Value = 33125

$ws.Range("E20").Value = "2001"
$ws.Range("F20").Value = 33125

# Update the totals: total "Valor Mora" and the period count.
$ws.Range("E11").Value = 144646
$ws.Range("F13").Value = 5
